# Omega Centauri Data.xlsx - progress update on two_body_utils.py / two_body.ipynb
# Reworks the "Assumptions" sheet (sheet2) with new columns of notes about
# the distance / projected-radius / velocity measurements used for the
# lower-limit central-mass estimate, and tidies up the workbook-level
# view state.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. Assumptions sheet (sheet2) column widths for the new note columns
# ---------------------------------------------------------------------
$ws2.Columns.Item(4).ColumnWidth = 55.85
$ws2.Columns.Item(5).ColumnWidth = 29.8

# ---------------------------------------------------------------------
# 2. Clear out the old single "Other" assumptions column (E2:E12) - its
#    content is being replaced by the new C/D/E layout below.
# ---------------------------------------------------------------------
$ws2.Range("E2:E12").ClearContents()

# ---------------------------------------------------------------------
# 3. New cell values
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "Assumptions"
$ws2.Range("A2").Value = "Old email thread"
$ws2.Range("C2").Value = "IMBH_Omega_Cen"
$ws2.Range("A3").Value = "binary barycenter velocity is zero"
$ws2.Range("C3").Value = "m_primary = 8200*M_SUN"
$ws2.Range("D3").Value = "lower limit esitmate for central mass"
$ws2.Range("A4").Value = "mass is 2e4 Msun "
$ws2.Range("C4").Value = "distance_kpc = 5.43 * u.kpc  "
$ws2.Range("D4").Value = "distance of omega centauri from earth"
$ws2.Range("C5").Value = 'r_proj_starA = 0.265"'
$ws2.Range("D5").Value = "angular distance from central mass for the two highest velocity stars "
$ws2.Range("E5").Value = 'Approximate average = 0.5"'
$ws2.Range("C6").Value = 'r_proj_starC = 0.870"'
$ws2.Range("C7").Value = "v_2D_starA = 113.0 km/s"
$ws2.Range("D7").Value = "2D velocities for the two highest velocity stars"
$ws2.Range("E7").Value = "Approximate average = 100 km/s"
$ws2.Range("C8").Value = "v_2D_starC = 94.9 km/s"

# ---------------------------------------------------------------------
# 4. Formatting - apply in this order so new style records are appended
#    to styles.xml in the same sequence the refreshed workbook uses.
# ---------------------------------------------------------------------

# vertical=top
$ws2.Range("D3").VerticalAlignment = -4160
$ws2.Range("A4").VerticalAlignment = -4160
$ws2.Range("D4").VerticalAlignment = -4160

# vertical=top + wrap
$ws2.Range("C4").VerticalAlignment = -4160
$ws2.Range("C4").WrapText = $true

# horizontal=left + vertical=top + wrap
$ws2.Range("D5:D6").HorizontalAlignment = -4131
$ws2.Range("D5:D6").VerticalAlignment = -4160
$ws2.Range("D5:D6").WrapText = $true

# bold + horizontal=center (title band)
$ws2.Range("A1:D1").Font.Bold = $true
$ws2.Range("A1:D1").HorizontalAlignment = -4108

# horizontal=left + vertical=top (no wrap)
$ws2.Range("E5:E6").HorizontalAlignment = -4131
$ws2.Range("E5:E6").VerticalAlignment = -4160
$ws2.Range("D7:D8").HorizontalAlignment = -4131
$ws2.Range("D7:D8").VerticalAlignment = -4160
$ws2.Range("E7:E8").HorizontalAlignment = -4131
$ws2.Range("E7:E8").VerticalAlignment = -4160

# ---------------------------------------------------------------------
# 5. Row heights
# ---------------------------------------------------------------------
$ws2.Rows.Item(4).RowHeight = 17.25

# ---------------------------------------------------------------------
# 6. Merge the grouped note cells
# ---------------------------------------------------------------------
$ws2.Range("D5:D6").Merge()
$ws2.Range("A1:D1").Merge()
$ws2.Range("D7:D8").Merge()
$ws2.Range("E5:E6").Merge()
$ws2.Range("E7:E8").Merge()

# ---------------------------------------------------------------------
# 7. Sheet view - scrolled down one row, selection parked on D19
# ---------------------------------------------------------------------
$ws2.Application.Goto($ws2.Range("A2"))
$ws2.Range("D19").Select()

# ---------------------------------------------------------------------
# 8. v2D and a_mag sheet (sheet1) is unchanged in content; re-touching
#    the header bands here is a no-op (same "horizontal=center" look).
# ---------------------------------------------------------------------
$ws1.Range("B1:C1").HorizontalAlignment = -4108
$ws1.Range("B11:G11").HorizontalAlignment = -4108
